# Auto-generated Excel COM-interop edit script
# Applies numeric cell updates per the authoritative diff, sheet by sheet.

$wb = $excel.ActiveWorkbook


$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 881.35895
$ws.Range("J17").Value = 866.8312
$ws.Range("L17").Value = 2600.4936
$ws.Range("N17").Value = -2936.4936
$ws.Range("H33").Value = 499.48486
$ws.Range("I33").Value = 555.1852
$ws.Range("K33").Value = 555.1852
$ws.Range("M33").Value = -326.1852
$ws.Range("H43").Value = 6189.6
$ws.Range("I43").Value = 20001
$ws.Range("K43").Value = 20001
$ws.Range("M43").Value = -19932
$ws.Range("H58").Value = 60
$ws.Range("I58").Value = 60
$ws.Range("K58").Value = 180
$ws.Range("M58").Value = -30
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("N65").ClearContents()
$ws.Range("H76").Value = 9998
$ws.Range("I76").Value = 9496.75
$ws.Range("K76").Value = 9496.75
$ws.Range("M76").Value = -9181.75
$ws.Range("H79").Value = 9998
$ws.Range("I79").Value = 9496.75
$ws.Range("K79").Value = 9496.75
$ws.Range("M79").Value = -8404.75
$ws.Range("H86").Value = 4525.4287
$ws.Range("I86").Value = 4556.25
$ws.Range("K86").Value = 4556.25
$ws.Range("M86").Value = -3433.25
$ws.Range("H89").Value = 4525.4287
$ws.Range("I89").Value = 4556.25
$ws.Range("K89").Value = 22781.25
$ws.Range("M89").Value = -17165.25
$ws.Range("H98").Value = 1391.4
$ws.Range("I98").Value = 1158.4166
$ws.Range("J98").Value = 2323.3333
$ws.Range("K98").Value = 1158.4166
$ws.Range("L98").Value = 2323.3333
$ws.Range("M98").Value = 339.5834
$ws.Range("N98").Value = -5319.3333
$ws.Range("H101").Value = 974.1429000000001
$ws.Range("I101").Value = 1086.6666
$ws.Range("K101").Value = 3259.9998
$ws.Range("M101").Value = -1637.9998
$ws.Range("H111").Value = 1835
$ws.Range("I111").Value = 1763.5
$ws.Range("J111").Value = 1882.6666
$ws.Range("K111").Value = 5290.5
$ws.Range("L111").Value = 5647.9998
$ws.Range("M111").Value = -2223.5
$ws.Range("N111").Value = -11781.9998
$ws.Range("H116").Value = 28500.334
$ws.Range("J116").Value = 31699.8
$ws.Range("L116").Value = 31699.8
$ws.Range("N116").Value = -38583.8
$ws.Range("H122").Value = 1391.4
$ws.Range("I122").Value = 1158.4166
$ws.Range("J122").Value = 2323.3333
$ws.Range("K122").Value = 3475.2498
$ws.Range("L122").Value = 6969.999899999999
$ws.Range("M122").Value = -1025.2498
$ws.Range("N122").Value = -11869.9999
$ws.Range("H127").Value = 1336.7646
$ws.Range("I127").Value = 980.3570999999999
$ws.Range("J127").Value = 3000
$ws.Range("K127").Value = 2941.0713
$ws.Range("L127").Value = 9000
$ws.Range("M127").Value = 2018.9287
$ws.Range("N127").Value = -18920
$ws.Range("H132").Value = 6744.83
$ws.Range("I132").Value = 4311.853
$ws.Range("J132").Value = 11098.579
$ws.Range("K132").Value = 12935.559
$ws.Range("L132").Value = 33295.737
$ws.Range("M132").Value = -10405.559
$ws.Range("N132").Value = -38355.737
$ws.Range("H138").Value = 796615.9399999999
$ws.Range("I138").Value = 2407.5
$ws.Range("J138").Value = 1013218.25
$ws.Range("K138").Value = 7222.5
$ws.Range("L138").Value = 3039654.75
$ws.Range("M138").Value = -2082.5
$ws.Range("N138").Value = -3049934.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6041.2236
$ws.Range("I32").Value = 4876.746
$ws.Range("K32").Value = 4876.746
$ws.Range("M32").Value = -4589.746
$ws.Range("H45").Value = 9074
$ws.Range("I45").Value = 10373.467
$ws.Range("K45").Value = 10373.467
$ws.Range("M45").Value = -9996.467000000001
$ws.Range("H61").Value = 3302.375
$ws.Range("I61").Value = 3239
$ws.Range("K61").Value = 3239
$ws.Range("M61").Value = -3027
$ws.Range("H74").Value = 2941.3572
$ws.Range("I74").Value = 2552.2307
$ws.Range("K74").Value = 2552.2307
$ws.Range("M74").Value = -1678.2307
$ws.Range("H77").Value = 2941.3572
$ws.Range("I77").Value = 2552.2307
$ws.Range("K77").Value = 12761.1535
$ws.Range("M77").Value = -8393.1535
$ws.Range("H110").Value = 3750
$ws.Range("I110").Value = 4875
$ws.Range("K110").Value = 4875
$ws.Range("M110").Value = -2830
$ws.Range("H122").Value = 30844.5
$ws.Range("I122").Value = 2786.5
$ws.Range("K122").Value = 8359.5
$ws.Range("M122").Value = -5909.5
$ws.Range("H132").Value = 2258.4783
$ws.Range("I132").Value = 2258.4783
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 6775.4349
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -4245.4349
$ws.Range("N132").ClearContents()
$ws.Range("H136").Value = 3302.375
$ws.Range("I136").Value = 3239
$ws.Range("K136").Value = 9717
$ws.Range("M136").Value = -7167
$ws.Range("H138").Value = 32499.5
$ws.Range("J138").Value = 32499.5
$ws.Range("L138").Value = 32499.5
$ws.Range("N138").Value = -42779.5
$ws.Range("H139").Value = 114972.875
$ws.Range("J139").Value = 114972.875
$ws.Range("L139").Value = 114972.875
$ws.Range("N139").Value = -125252.875

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H25").Value = 8589
$ws.Range("I25").Value = 8589
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 8589
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = -8354
$ws.Range("N25").ClearContents()
$ws.Range("H63").Value = 60000
$ws.Range("J63").Value = 60000
$ws.Range("L63").Value = 60000
$ws.Range("N63").Value = -61372
$ws.Range("H66").Value = 60000
$ws.Range("J66").Value = 60000
$ws.Range("L66").Value = 180000
$ws.Range("N66").Value = -186864
$ws.Range("H86").Value = 2192.1333
$ws.Range("I86").Value = 2017.7273
$ws.Range("K86").Value = 2017.7273
$ws.Range("M86").Value = -894.7273
$ws.Range("H89").Value = 2192.1333
$ws.Range("I89").Value = 2017.7273
$ws.Range("K89").Value = 10088.6365
$ws.Range("M89").Value = -4472.636500000001
$ws.Range("H99").Value = 1533.4445
$ws.Range("I99").Value = 1325.4642
$ws.Range("K99").Value = 1325.4642
$ws.Range("M99").Value = 172.5358000000001
$ws.Range("H112").Value = 175000
$ws.Range("J112").Value = 175000
$ws.Range("L112").Value = 175000
$ws.Range("N112").Value = -177954
$ws.Range("H134").Value = 6394
$ws.Range("I134").Value = 6192
$ws.Range("J134").Value = 7336.6665
$ws.Range("K134").Value = 18576
$ws.Range("L134").Value = 22009.9995
$ws.Range("M134").Value = -16041
$ws.Range("N134").Value = -27079.9995

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 5166
$ws.Range("I4").Value = 3500
$ws.Range("K4").Value = 3500
$ws.Range("M4").Value = -3388
$ws.Range("H16").Value = 2406.4546
$ws.Range("I16").Value = 1719
$ws.Range("J16").Value = 5500
$ws.Range("K16").Value = 1719
$ws.Range("L16").Value = 5500
$ws.Range("M16").Value = -1432
$ws.Range("N16").Value = -6074
$ws.Range("H22").Value = 215.46666
$ws.Range("I22").Value = 231.72728
$ws.Range("J22").Value = 170.75
$ws.Range("K22").Value = 231.72728
$ws.Range("L22").Value = 170.75
$ws.Range("M22").Value = 118.27272
$ws.Range("N22").Value = -870.75
$ws.Range("H31").Value = 1965.5667
$ws.Range("I31").Value = 1396.2
$ws.Range("K31").Value = 1396.2
$ws.Range("M31").Value = -1101.2
$ws.Range("H34").Value = 1965.5667
$ws.Range("I34").Value = 1396.2
$ws.Range("K34").Value = 1396.2
$ws.Range("M34").Value = -1194.2
$ws.Range("H35").Value = 816
$ws.Range("I35").Value = 543
$ws.Range("K35").Value = 543
$ws.Range("M35").Value = -249
$ws.Range("H58").Value = 2223.2942
$ws.Range("J58").Value = 2178.4285
$ws.Range("L58").Value = 2178.4285
$ws.Range("N58").Value = -2584.4285
$ws.Range("H86").Value = 3706253.8
$ws.Range("I86").Value = 5130673
$ws.Range("K86").Value = 5130673
$ws.Range("M86").Value = -5129550
$ws.Range("H89").Value = 3706253.8
$ws.Range("I89").Value = 5130673
$ws.Range("K89").Value = 25653365
$ws.Range("M89").Value = -25647749
$ws.Range("H94").Value = 4275.231
$ws.Range("J94").Value = 4584.75
$ws.Range("L94").Value = 4584.75
$ws.Range("N94").Value = -5486.75
$ws.Range("H99").Value = 11658.314
$ws.Range("I99").Value = 6475.1665
$ws.Range("K99").Value = 6475.1665
$ws.Range("M99").Value = -4977.1665
$ws.Range("H111").Value = 193318
$ws.Range("J111").Value = 193318
$ws.Range("L111").Value = 193318
$ws.Range("N111").Value = -201498
$ws.Range("H113").Value = 2406.4546
$ws.Range("I113").Value = 1719
$ws.Range("J113").Value = 5500
$ws.Range("K113").Value = 1719
$ws.Range("L113").Value = 5500
$ws.Range("M113").Value = 451
$ws.Range("N113").Value = -9840
$ws.Range("H126").Value = 11658.314
$ws.Range("I126").Value = 6475.1665
$ws.Range("K126").Value = 19425.4995
$ws.Range("M126").Value = -16955.4995
$ws.Range("H131").Value = 75050.664
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 75050.664
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 75050.664
$ws.Range("M131").ClearContents()
$ws.Range("N131").Value = -85130.664
$ws.Range("H132").Value = 6467.4
$ws.Range("I132").Value = 3476.2104
$ws.Range("J132").Value = 13153.588
$ws.Range("K132").Value = 10428.6312
$ws.Range("L132").Value = 39460.764
$ws.Range("M132").Value = -7898.6312
$ws.Range("N132").Value = -44520.764
$ws.Range("H134").Value = 1926.4694
$ws.Range("I134").Value = 1793.1428
$ws.Range("K134").Value = 5379.428400000001
$ws.Range("M134").Value = -2844.428400000001
$ws.Range("H136").Value = 2223.2942
$ws.Range("J136").Value = 2178.4285
$ws.Range("L136").Value = 6535.2855
$ws.Range("N136").Value = -11635.2855

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 10000
$ws.Range("J75").Value = 10000
$ws.Range("L75").Value = 30000
$ws.Range("N75").Value = -31996
$ws.Range("H78").Value = 10000
$ws.Range("J78").Value = 10000
$ws.Range("L78").Value = 90000
$ws.Range("N78").Value = -99984
$ws.Range("H81").Value = 1566
$ws.Range("I81").Value = 1566
$ws.Range("K81").Value = 4698
$ws.Range("M81").Value = -3575
$ws.Range("H84").Value = 1566
$ws.Range("I84").Value = 1566
$ws.Range("K84").Value = 14094
$ws.Range("M84").Value = -8478
$ws.Range("H98").Value = 920.25
$ws.Range("J98").Value = 896
$ws.Range("L98").Value = 2688
$ws.Range("N98").Value = -5684
$ws.Range("H113").Value = 1129
$ws.Range("I113").Value = 738.8570999999999
$ws.Range("J113").Value = 1519.1428
$ws.Range("K113").Value = 2216.5713
$ws.Range("L113").Value = 4557.428400000001
$ws.Range("M113").Value = -46.57129999999961
$ws.Range("N113").Value = -8897.428400000001
$ws.Range("H141").Value = 13291.9
$ws.Range("J141").Value = 5999.5
$ws.Range("L141").Value = 17998.5
$ws.Range("N141").Value = -28358.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H6").Value = 2949
$ws.Range("I6").Value = 3999
$ws.Range("J6").Value = 1899
$ws.Range("K6").Value = 3999
$ws.Range("L6").Value = 1899
$ws.Range("M6").Value = -3886
$ws.Range("N6").Value = -2125
$ws.Range("H11").Value = 10554600
$ws.Range("I11").Value = 5104400
$ws.Range("K11").Value = 5104400
$ws.Range("M11").Value = -5104261
$ws.Range("H12").Value = 167.33333
$ws.Range("I12").Value = 3
$ws.Range("J12").Value = 249.5
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 249.5
$ws.Range("M12").Value = 137
$ws.Range("N12").Value = -529.5
$ws.Range("H16").Value = 2949
$ws.Range("I16").Value = 3999
$ws.Range("J16").Value = 1899
$ws.Range("K16").Value = 3999
$ws.Range("L16").Value = 1899
$ws.Range("M16").Value = -3749
$ws.Range("N16").Value = -2399
$ws.Range("H20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").ClearContents()
$ws.Range("H24").Value = 18856.857
$ws.Range("I24").Value = 14000
$ws.Range("J24").Value = 19666.334
$ws.Range("K24").Value = 14000
$ws.Range("L24").Value = 19666.334
$ws.Range("M24").Value = -13827
$ws.Range("N24").Value = -20012.334
$ws.Range("H31").Value = 917.25
$ws.Range("I31").Value = 917.25
$ws.Range("K31").Value = 917.25
$ws.Range("M31").Value = -625.25
$ws.Range("H37").Value = 917.25
$ws.Range("I37").Value = 917.25
$ws.Range("K37").Value = 917.25
$ws.Range("M37").Value = -640.25
$ws.Range("H51").Value = 175000
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 175000
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 175000
$ws.Range("M51").ClearContents()
$ws.Range("N51").Value = -176018
$ws.Range("H70").Value = 5116.6484
$ws.Range("I70").Value = 5063.8335
$ws.Range("K70").Value = 5063.8335
$ws.Range("M70").Value = -4793.8335
$ws.Range("H73").Value = 5116.6484
$ws.Range("I73").Value = 5063.8335
$ws.Range("K73").Value = 5063.8335
$ws.Range("M73").Value = -4127.8335
$ws.Range("H113").Value = 2926.2
$ws.Range("I113").Value = 2686.5557
$ws.Range("K113").Value = 2686.5557
$ws.Range("M113").Value = -516.5556999999999
$ws.Range("H122").Value = 1249
$ws.Range("I122").Value = 1249
$ws.Range("K122").Value = 3747
$ws.Range("M122").Value = -1297
$ws.Range("H126").Value = 5546.25
$ws.Range("J126").Value = 4970.3335
$ws.Range("L126").Value = 14911.0005
$ws.Range("N126").Value = -19851.0005
$ws.Range("H132").Value = 5256.615
$ws.Range("I132").Value = 2973.1555
$ws.Range("J132").Value = 19936
$ws.Range("K132").Value = 8919.466499999999
$ws.Range("L132").Value = 59808
$ws.Range("M132").Value = -6389.466499999999
$ws.Range("N132").Value = -64868
$ws.Range("H133").Value = 100000
$ws.Range("J133").Value = 100000
$ws.Range("L133").Value = 100000
$ws.Range("N133").Value = -110120

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 8996
$ws.Range("I7").Value = 8996
$ws.Range("K7").Value = 8996
$ws.Range("M7").Value = -8884
$ws.Range("H23").Value = 18001.5
$ws.Range("I23").Value = 13002
$ws.Range("J23").Value = 33000
$ws.Range("K23").Value = 13002
$ws.Range("L23").Value = 33000
$ws.Range("M23").Value = -12772
$ws.Range("N23").Value = -33460
$ws.Range("H40").Value = 6499.5
$ws.Range("I40").Value = 6499.5
$ws.Range("K40").Value = 6499.5
$ws.Range("M40").Value = -6363.5
$ws.Range("H43").Value = 20000
$ws.Range("J43").Value = 20000
$ws.Range("L43").Value = 20000
$ws.Range("N43").Value = -20386
$ws.Range("H122").Value = 8756.714
$ws.Range("I122").Value = 9416.166999999999
$ws.Range("J122").Value = 4800
$ws.Range("K122").Value = 28248.501
$ws.Range("L122").Value = 14400
$ws.Range("M122").Value = -25798.501
$ws.Range("N122").Value = -19300
$ws.Range("H126").Value = 8996
$ws.Range("I126").Value = 8996
$ws.Range("K126").Value = 26988
$ws.Range("M126").Value = -24518
$ws.Range("H136").Value = 3605501.8
$ws.Range("I136").Value = 5459258.5
$ws.Range("K136").Value = 16377775.5
$ws.Range("M136").Value = -16375225.5
$ws.Range("H137").Value = 54524.25
$ws.Range("J137").Value = 59365.668
$ws.Range("L137").Value = 59365.668
$ws.Range("N137").Value = -69565.66800000001
$ws.Range("H138").Value = 147666
$ws.Range("J138").Value = 147666
$ws.Range("L138").Value = 147666
$ws.Range("N138").Value = -157946

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H37").Value = 49999
$ws.Range("I37").Value = 49999
$ws.Range("K37").Value = 49999
$ws.Range("M37").Value = -49796
$ws.Range("H75").Value = 18950
$ws.Range("I75").Value = 18000
$ws.Range("K75").Value = 18000
$ws.Range("M75").Value = -17064
$ws.Range("H78").Value = 18950
$ws.Range("I78").Value = 18000
$ws.Range("K78").Value = 54000
$ws.Range("M78").Value = -49320
$ws.Range("H122").Value = 37542.43
$ws.Range("I122").Value = 35666.668
$ws.Range("J122").Value = 38949.25
$ws.Range("K122").Value = 107000.004
$ws.Range("L122").Value = 116847.75
$ws.Range("M122").Value = -104550.004
$ws.Range("N122").Value = -121747.75
$ws.Range("H126").Value = 4660.7646
$ws.Range("J126").Value = 3301.7144
$ws.Range("L126").Value = 9905.143199999999
$ws.Range("N126").Value = -14845.1432
$ws.Range("H132").Value = 5401.6523
$ws.Range("I132").Value = 4090.7222
$ws.Range("K132").Value = 12272.1666
$ws.Range("M132").Value = -9742.1666
$ws.Range("H136").Value = 1440.5964
$ws.Range("I136").Value = 1643.7556
$ws.Range("K136").Value = 4931.266799999999
$ws.Range("M136").Value = -2381.266799999999
$ws.Range("H137").Value = 122214.5
$ws.Range("J137").Value = 122214.5
$ws.Range("L137").Value = 122214.5
$ws.Range("N137").Value = -132414.5
$ws.Range("H139").Value = 54698.6
$ws.Range("J139").Value = 54698.6
$ws.Range("L139").Value = 54698.6
$ws.Range("N139").Value = -64978.6
